$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1035.591
$ws.Range("I41").Value = 1743.625
$ws.Range("J41").Value = 631
$ws.Range("K41").Value = 1743.625
$ws.Range("L41").Value = 631
$ws.Range("M41").Value = -1303.625
$ws.Range("N41").Value = -1511
$ws.Range("H92").Value = 548.52
$ws.Range("I92").Value = 547.94446
$ws.Range("K92").Value = 547.94446
$ws.Range("M92").Value = 700.05554
$ws.Range("H132").Value = 5958609.5
$ws.Range("I132").Value = 6950794.5
$ws.Range("K132").Value = 20852383.5
$ws.Range("M132").Value = -20849853.5
$ws.Range("H138").Value = 2813.9849
$ws.Range("J138").Value = 3802.2368
$ws.Range("L138").Value = 11406.7104
$ws.Range("N138").Value = -21686.7104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22619.414
$ws.Range("I32").Value = 6377.239
$ws.Range("K32").Value = 6377.239
$ws.Range("M32").Value = -6090.239
$ws.Range("H61").Value = 1852.7894
$ws.Range("I61").Value = 1471.6666
$ws.Range("K61").Value = 1471.6666
$ws.Range("M61").Value = -1259.6666
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("H122").Value = 2535.0476
$ws.Range("I122").Value = 2118.4
$ws.Range("J122").Value = 3576.6667
$ws.Range("K122").Value = 6355.200000000001
$ws.Range("L122").Value = 10730.0001
$ws.Range("M122").Value = -3905.200000000001
$ws.Range("N122").Value = -15630.0001
$ws.Range("H136").Value = 1852.7894
$ws.Range("I136").Value = 1471.6666
$ws.Range("K136").Value = 4414.9998
$ws.Range("M136").Value = -1864.9998
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("H86").Value = 56840.3
$ws.Range("I86").Value = 70519.125
$ws.Range("K86").Value = 70519.125
$ws.Range("M86").Value = -69396.125
$ws.Range("H89").Value = 56840.3
$ws.Range("I89").Value = 70519.125
$ws.Range("K89").Value = 352595.625
$ws.Range("M89").Value = -346979.625
$ws.Range("H94").Value = 841.8333
$ws.Range("I94").Value = 841.73334
$ws.Range("J94").Value = 842.3333
$ws.Range("K94").Value = 841.73334
$ws.Range("L94").Value = 842.3333
$ws.Range("M94").Value = -390.73334
$ws.Range("N94").Value = -1744.3333
$ws.Range("H102").Value = 15018.667
$ws.Range("I102").Value = 15018.667
$ws.Range("K102").Value = 15018.667
$ws.Range("M102").Value = -11773.667
$ws.Range("M16").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 11999.2
$ws.Range("J70").Value = 11999.2
$ws.Range("L70").Value = 11999.2
$ws.Range("N70").Value = -12629.2
$ws.Range("H73").Value = 11999.2
$ws.Range("J73").Value = 11999.2
$ws.Range("L73").Value = 11999.2
$ws.Range("N73").Value = -14183.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 85500.164
$ws.Range("J70").Value = 2914.2856
$ws.Range("L70").Value = 8742.856800000001
$ws.Range("N70").Value = -9372.856800000001
$ws.Range("H73").Value = 85500.164
$ws.Range("J73").Value = 2914.2856
$ws.Range("L73").Value = 8742.856800000001
$ws.Range("N73").Value = -10926.8568
$ws.Range("H80").Value = 8041.7334
$ws.Range("J80").Value = 9802.166999999999
$ws.Range("L80").Value = 29406.501
$ws.Range("N80").Value = -31278.501
$ws.Range("H83").Value = 8041.7334
$ws.Range("J83").Value = 9802.166999999999
$ws.Range("L83").Value = 88219.503
$ws.Range("N83").Value = -97579.503
$ws.Range("H122").Value = 654.1429000000001
$ws.Range("I122").Value = 594.75
$ws.Range("K122").Value = 5352.75
$ws.Range("M122").Value = -2902.75
$ws.Range("H137").Value = 5318402.5
$ws.Range("I137").Value = 68582.664
$ws.Range("J137").Value = 25005228
$ws.Range("K137").Value = 205747.992
$ws.Range("L137").Value = 75015684
$ws.Range("M137").Value = -200647.992
$ws.Range("N137").Value = -75025884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 28000
$ws.Range("J32").Value = 28000
$ws.Range("L32").Value = 28000
$ws.Range("N32").Value = -28592
$ws.Range("H42").Value = 38540.75
$ws.Range("J42").Value = 38540.75
$ws.Range("L42").Value = 38540.75
$ws.Range("N42").Value = -39510.75
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
$ws.Range("H115").Value = 38540.75
$ws.Range("J115").Value = 38540.75
$ws.Range("L115").Value = 38540.75
$ws.Range("N115").Value = -40890.75
$ws.Range("H132").Value = 2866.9524
$ws.Range("I132").Value = 2132.1035
$ws.Range("K132").Value = 6396.310500000001
$ws.Range("M132").Value = -3866.310500000001
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2743.1765
$ws.Range("I132").Value = 2796
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 8388
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -5858
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 13700
$ws.Range("J25").Value = 13700
$ws.Range("L25").Value = 13700
$ws.Range("N25").Value = -14286
$ws.Range("H27").Value = 35649.5
$ws.Range("J27").Value = 35649.5
$ws.Range("L27").Value = 35649.5
$ws.Range("N27").Value = -35787.5
$ws.Range("H69").Value = 24196
$ws.Range("J69").Value = 24196
$ws.Range("L69").Value = 24196
$ws.Range("N69").Value = -25694
$ws.Range("H72").Value = 24196
$ws.Range("J72").Value = 24196
$ws.Range("L72").Value = 72588
$ws.Range("N72").Value = -80076
$ws.Range("H107").Value = 55971.832
$ws.Range("I107").Value = 406.2
$ws.Range("J107").Value = 333800
$ws.Range("K107").Value = 1218.6
$ws.Range("L107").Value = 1001400
$ws.Range("M107").Value = 701.4000000000001
$ws.Range("N107").Value = -1005240
$ws.Range("H113").Value = 590.3125
$ws.Range("I113").Value = 435.63635
$ws.Range("J113").Value = 930.6
$ws.Range("K113").Value = 1306.90905
$ws.Range("L113").Value = 2791.8
$ws.Range("M113").Value = 863.09095
$ws.Range("N113").Value = -7131.8
$ws.Range("H115").Value = 36332.582
$ws.Range("J115").Value = 36332.582
$ws.Range("L115").Value = 36332.582
$ws.Range("N115").Value = -39466.582

Write-Host "Applied all changes"